$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7251761691223575
$ws.Range("D2").Value = "Below Median"
$ws.Range("C3").Value = 0.9199231262011531
$ws.Range("D3").Value = "Below Median"
$ws.Range("C4").Value = 2.65829596412556
$ws.Range("D4").Value = "1st Tier"
$ws.Range("C5").Value = 2.026905829596413
$ws.Range("D5").Value = "2nd Tier"
$ws.Range("C6").Value = 3.136023916292975
$ws.Range("D6").Value = "1st Tier"
$ws.Range("C7").Value = 2.429212043561819
$ws.Range("D7").Value = "1st Tier"
$ws.Range("C8").Value = 1.596412556053812
$ws.Range("D8").Value = "3rd Tier"
$ws.Range("C9").Value = 4.340807174887892
$ws.Range("D9").Value = "1st Tier"
$ws.Range("C10").Value = 0.6278026905829597
$ws.Range("D10").Value = "Below Median"
$ws.Range("C11").Value = 0.7623318385650224
$ws.Range("D11").Value = "Below Median"
$ws.Range("C12").Value = 0.8456117873158232
$ws.Range("D12").Value = "Below Median"
$ws.Range("C13").Value = 0.9820627802690582
$ws.Range("D13").Value = "Below Median"
$ws.Range("C14").Value = 0.9820627802690582
$ws.Range("D14").Value = "Below Median"
$ws.Range("C15").Value = 0.8998505231689088
$ws.Range("D15").Value = "Below Median"
$ws.Range("C16").Value = 1.086995515695067
$ws.Range("D16").Value = "4th Tier"
$ws.Range("C17").Value = 1.051569506726457
$ws.Range("D17").Value = "4th Tier"
$ws.Range("C18").Value = 1.237668161434978
$ws.Range("D18").Value = "4th Tier"
$ws.Range("C19").Value = 0.9108121574489287
$ws.Range("D19").Value = "Below Median"
$ws.Range("C20").Value = 0.7593423019431988
$ws.Range("D20").Value = "Below Median"
$ws.Range("C21").Value = 0.7892376681614349
$ws.Range("D21").Value = "Below Median"
$ws.Range("C22").Value = 0.8891736066623959
$ws.Range("D22").Value = "Below Median"
$ws.Range("C23").Value = 0.8379244074311339
$ws.Range("D23").Value = "Below Median"
$ws.Range("C24").Value = 2.895067264573991
$ws.Range("D24").Value = "1st Tier"
$ws.Range("C25").Value = 1.796284433055733
$ws.Range("D25").Value = "2nd Tier"
$ws.Range("C26").Value = 3.019431988041854
$ws.Range("D26").Value = "1st Tier"
$ws.Range("C27").Value = 1.461883408071749
$ws.Range("D27").Value = "3rd Tier"
$ws.Range("C28").Value = 2.125560538116592
$ws.Range("D28").Value = "2nd Tier"
$ws.Range("C29").Value = 0.9650224215246637
$ws.Range("D29").Value = "Below Median"
$ws.Range("C30").Value = 2.089686098654708
$ws.Range("D30").Value = "2nd Tier"
$ws.Range("C31").Value = 1.285500747384156
$ws.Range("D31").Value = "3rd Tier"
$ws.Range("C32").Value = 2.234977578475336
$ws.Range("D32").Value = "2nd Tier"
$ws.Range("C33").Value = 0.5769805680119581
$ws.Range("D33").Value = "Below Median"
$ws.Range("C34").Value = 2.078155028827675
$ws.Range("D34").Value = "2nd Tier"
$ws.Range("C35").Value = 0.726457399103139
$ws.Range("D35").Value = "Below Median"
$ws.Range("C36").Value = 2.208840486867393
$ws.Range("D36").Value = "2nd Tier"
$ws.Range("C37").Value = 0.5739910313901345
$ws.Range("D37").Value = "Below Median"
$ws.Range("C38").Value = 1.409352978859705
$ws.Range("D38").Value = "3rd Tier"
$ws.Range("C39").Value = 1.194106342088405
$ws.Range("D39").Value = "4th Tier"
$ws.Range("C40").Value = 2.787443946188341
$ws.Range("D40").Value = "1st Tier"
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = "4th Tier"
$ws.Range("C42").Value = 1.1898355754858
$ws.Range("D42").Value = "4th Tier"
$ws.Range("C43").Value = 1.809417040358744
$ws.Range("D43").Value = "2nd Tier"
$ws.Range("C44").Value = 1.456176110884631
$ws.Range("D44").Value = "3rd Tier"
$ws.Range("C45").Value = 1.069058295964125
$ws.Range("D45").Value = "4th Tier"
$ws.Range("C46").Value = 0.2654708520179372
$ws.Range("D46").Value = "Below Median"
$ws.Range("C47").Value = 0.7294469357249626
$ws.Range("D47").Value = "Below Median"
$ws.Range("C48").Value = 1.257847533632287
$ws.Range("D48").Value = "3rd Tier"
$ws.Range("C49").Value = 0.6744394618834081
$ws.Range("D49").Value = "Below Median"
$ws.Range("C50").Value = 0.9125560538116592
$ws.Range("D50").Value = "Below Median"
$ws.Range("C51").Value = 0.2750373692077728
$ws.Range("D51").Value = "Below Median"
$ws.Range("C52").Value = 0.8968609865470852
$ws.Range("D52").Value = "Below Median"
$ws.Range("C53").Value = 2.517189835575486
$ws.Range("D53").Value = "1st Tier"
$ws.Range("C54").Value = 0.2331838565022422
$ws.Range("D54").Value = "Below Median"
$ws.Range("C55").Value = 1.522101217168482
$ws.Range("D55").Value = "3rd Tier"
$ws.Range("C56").Value = 0.9147982062780269
$ws.Range("D56").Value = "Below Median"
$ws.Range("C57").Value = 1.280717488789238
$ws.Range("D57").Value = "3rd Tier"
$ws.Range("C58").Value = 0.5704035874439461
$ws.Range("D58").Value = "Below Median"
$ws.Range("C59").Value = 0.2600896860986547
$ws.Range("D59").Value = "Below Median"
$ws.Range("C60").Value = 0.5944907110826394
$ws.Range("D60").Value = "Below Median"
$ws.Range("C61").Value = 1.959641255605381
$ws.Range("D61").Value = "2nd Tier"
$ws.Range("C62").Value = 2.946188340807175
$ws.Range("D62").Value = "1st Tier"
$ws.Range("C63").Value = 0.8143497757847533
$ws.Range("D63").Value = "Below Median"
$ws.Range("C64").Value = 0.7461883408071749
$ws.Range("D64").Value = "Below Median"
$ws.Range("C65").Value = 0.6816143497757847
$ws.Range("D65").Value = "Below Median"
$ws.Range("C66").Value = 1.174887892376682
$ws.Range("D66").Value = "4th Tier"
$ws.Range("C67").Value = 1.818834080717489
$ws.Range("D67").Value = "2nd Tier"
$ws.Range("C68").Value = 0.8789237668161435
$ws.Range("D68").Value = "Below Median"
$ws.Range("C69").Value = 0.2798206278026906
$ws.Range("D69").Value = "Below Median"
$ws.Range("C70").Value = 0.7713004484304933
$ws.Range("D70").Value = "Below Median"
$ws.Range("C71").Value = 0.5847533632286995
$ws.Range("D71").Value = "Below Median"
$ws.Range("C72").Value = 1.356053811659193
$ws.Range("D72").Value = "3rd Tier"
$ws.Range("C73").Value = 0.75695067264574
$ws.Range("D73").Value = "Below Median"
$ws.Range("C74").Value = 1.191543882126842
$ws.Range("D74").Value = "4th Tier"
$ws.Range("C75").Value = 4.573991031390134
$ws.Range("D75").Value = "1st Tier"
$ws.Range("C76").Value = 0.8938714499252616
$ws.Range("D76").Value = "Below Median"
$ws.Range("C77").Value = 0.8968609865470852
$ws.Range("D77").Value = "Below Median"
$ws.Range("C78").Value = 1.704035874439462
$ws.Range("D78").Value = "3rd Tier"
$ws.Range("C79").Value = 0.57847533632287
$ws.Range("D79").Value = "Below Median"
$ws.Range("C80").Value = 1.174887892376682
$ws.Range("D80").Value = "4th Tier"
